$wb = $excel.ActiveWorkbook

# --- Sheet 1: OUTP1M_RATIO ---
$ws1 = $wb.Worksheets.Item("OUTP1M_RATIO")
$values1 = @(
    3.17628802895906,
    2.78638161887961,
    2.1042176646398,
    1.69535510502432,
    1.10879308735565,
    0.995629744266522,
    0.869729180034146,
    0.78740217344344,
    0.696219201282135,
    0.598949283971256,
    0.471961108424601,
    0.324712925902315,
    0.141797149123868,
    0.11192129353216,
    0.109476867300401,
    0.0211655678607131
)
for ($i = 0; $i -lt $values1.Length; $i++) {
    $row = $i + 1
    $ws1.Cells.Item($row, 1).Value = $values1[$i]
}

# --- Sheet 2: CHRONIC_RATIO ---
$ws2 = $wb.Worksheets.Item("CHRONIC_RATIO")
$values2 = @(
    3.70257288893534,
    2.592355898322,
    1.89106947479405,
    1.44190343749358,
    1.19607128962009,
    1.03679147647512,
    0.909115267738115,
    0.805488667484156,
    0.739457764253924,
    0.68140363598806,
    0.60540587731287,
    0.469342523884538,
    0.378197546636995,
    0.29416551317241,
    0.130622340446943,
    0.118846518958624,
    0.00718987848318083
)
for ($i = 0; $i -lt $values2.Length; $i++) {
    $row = $i + 1
    $ws2.Cells.Item($row, 1).Value = $values2[$i]
}
